$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (FP1) - D2 gets the quoted "son 'chaud'" wording, E2 keeps wrap style
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "2x3W, son 'chaud' et non saturé"

# ---------------------------------------------------------------------------
# Row 3 (FC1) - new "latence" (écran tactile/boutons) function block
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Traiter les informations provenant de l'écran tactile et des boutons"
$ws.Range("C3").Value = "Faible latence et fluidité des contrôles`nMinimalisme des commandes possibles"
$ws.Range("D3").Value = "Moins de 100 ms entre l'appui et la réponse"
$ws.Range("E3").Value = "F1`nLatence : ±20ms"

$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $true
$ws.Range("E3").WrapText = $true

# ---------------------------------------------------------------------------
# Row 4 (FC2) - unchanged function, only moved here from the old row 3
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "S'intégrer au boitier en n'altérant pas le design de l'objet"

# ---------------------------------------------------------------------------
# Row 5 (FC3) - new "application compagnon" function block
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Etablir une connectivité et traiter les informations provenant de l'application compagnon"
$ws.Range("C5").Value = "Faible latence et fluidité des contrôles`nPortée importante`nMinimalisme des commandes possibles"
$ws.Range("D5").Value = "Moins de 100 ms entre l'appui et la réponse`nPortée de 40m en indoor"
$ws.Range("E5").Value = "F1`nLatence : ±20ms`nPortée : ±15m"

$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("C5").WrapText = $true
$ws.Range("D5").WrapText = $true
$ws.Range("E5").WrapText = $true

# ---------------------------------------------------------------------------
# Row 6 (FC4) - new "normes" function block
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "Respect des différentes normes en vigueur*"
$ws.Range("C6").Value = "Respect des normes"
$ws.Range("D6").Value = "Toutes les normes doivent être respectées"
$ws.Range("E6").Value = "F0"
$ws.Range("E6").WrapText = $true

# ---------------------------------------------------------------------------
# Row 7 (FC5) - new "tuner" function block
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Traiter les données provenant du tuner"
$ws.Range("C7").Value = "Traitement rapide et sans perte"
$ws.Range("D7").Value = "Pas de pertes de données"
$ws.Range("E7").Value = "F0"

$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("C7").WrapText = $true
$ws.Range("E7").WrapText = $true

# ---------------------------------------------------------------------------
# Row 9 - footnote about the applicable norms
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "*"
$ws.Range("B9").Value = "Normes CE principalement : `nCompatibilité électromagnétique (CEM) - 2014/30/UE`nÉquipements terminaux de télécommunication - 1999/5/CE "
$ws.Range("B9").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights for the newly wrapped rows
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# ---------------------------------------------------------------------------
# Column widths (B/C/D got wider to fit the new text)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 80.5
$ws.Columns.Item(3).ColumnWidth = 35.833333
$ws.Columns.Item(4).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Selection ends on B9, matching the final cursor position
# ---------------------------------------------------------------------------
[void]$ws.Range("B9").Select()
